$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of Box Office data (rank 10 - Big Hero 6)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Big Hero 6"
$ws.Range("C11").Value = "BV"
$ws.Range("D11").Value = 222527828
$ws.Range("E11").Value = 3773
$ws.Range("F11").Value = 56215889
$ws.Range("G11").Value = 3761
$ws.Range("H11").Value = 43411
$ws.Range("I11").Value = 43248

# Match number formats with the row above so the new cells render the same way
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E11").NumberFormat = $ws.Range("E10").NumberFormat
$ws.Range("F11").NumberFormat = $ws.Range("F10").NumberFormat
$ws.Range("G11").NumberFormat = $ws.Range("G10").NumberFormat
$ws.Range("H11").NumberFormat = $ws.Range("H10").NumberFormat
$ws.Range("I11").NumberFormat = $ws.Range("I10").NumberFormat

# Move the active selection to match the saved view
[void]$ws.Range("C21").Select()
